$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C3").Value = 17.530752052259572
$ws.Range("C5").Value = -0.5882549411371829
$ws.Range("C7").Value = 41.546788901173656
$ws.Range("C9").Value = -15.09018482792024
$ws.Range("C13").Value = 16.599017050968996
$ws.Range("C15").Value = -0.8280762962347101
$ws.Range("C17").Value = 17.645496456545658
$ws.Range("C19").Value = -21.242191927270767
$ws.Range("C23").Value = 16.599017050968996
$ws.Range("C25").Value = -0.8280762962347101
$ws.Range("C27").Value = 17.645496456545658
$ws.Range("C29").Value = -21.242191927270767
$ws.Range("C33").Value = 16.599017050968996
$ws.Range("C35").Value = -0.8280762962347101
$ws.Range("C37").Value = 17.645496456545658
$ws.Range("C39").Value = -21.242191927270767
$ws.Range("C43").Value = 17.132418618993512
$ws.Range("C45").Value = -0.6052579094319128
$ws.Range("C47").Value = 31.328557466508812
$ws.Range("C49").Value = -15.526352748065095
$ws.Range("C53").Value = 17.02389411303927
$ws.Range("C55").Value = -0.6412881530154779
$ws.Range("C57").Value = 28.544637198805777
$ws.Range("C59").Value = -16.450617037320242
$ws.Range("C62").Value = 11.411225144977507
$ws.Range("C63").Value = 28.544637198805777
$ws.Range("C64").Value = 35.15015499138343
$ws.Range("C69").Value = 49945.16433399591
$ws.Range("C70").Value = 2977293.533145299
$ws.Range("C71").Value = 2927348.3688113037
$ws.Range("C76").Value = 12889.249455863679

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C15").Value = -20.52196593367998
$ws.Range("C16").Value = -20.52196593367998
$ws.Range("C19").Value = -0.7999999999999999
$ws.Range("C20").Value = -0.7999999999999999

$ws = $wb.Worksheets.Item("FUEL TANK")
$ws.Range("C15").Value = -20.52196593367998
$ws.Range("C16").Value = -20.52196593367998
$ws.Range("C19").Value = -0.7999999999999999
$ws.Range("C20").Value = -0.7999999999999998

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = 12.95357995250832
$ws.Range("C6").Value = 12.953579952508276
$ws.Range("C7").Value = 16.416113852739358
$ws.Range("C8").Value = 16.416113852739354
$ws.Range("C9").Value = 16.416113852739354
$ws.Range("C10").Value = 16.41611385273935
$ws.Range("C23").Value = 16.416113852739358
